$d = $word.ActiveDocument

$pairs = @(
    @{old="14×69=966"; new="88×67=5896"},
    @{old="38×36=1368"; new="21×91=1911"},
    @{old="50×67=3350"; new="90×14=1260"},
    @{old="31×49=1519"; new="62×24=1488"},
    @{old="84×23=1932"; new="65×27=1755"},
    @{old="86×26=2236"; new="19×48=912"},
    @{old="30×55=1650"; new="61×30=1830"},
    @{old="96×97=9312"; new="38×80=3040"},
    @{old="75×81=6075"; new="53×81=4293"},
    @{old="75×59=4425"; new="78×60=4680"},
    @{old="23×70=1610"; new="69×36=2484"},
    @{old="49×24=1176"; new="48×36=1728"},
    @{old="77×95=7315"; new="59×79=4661"},
    @{old="58×46=2668"; new="34×90=3060"},
    @{old="93×64=5952"; new="14×88=1232"},
    @{old="59×85=5015"; new="23×38=874"},
    @{old="13×12=156"; new="39×99=3861"},
    @{old="50×98=4900"; new="62×48=2976"},
    @{old="20×11=220"; new="12×84=1008"},
    @{old="39×57=2223"; new="39×38=1482"},
    @{old="59×50=2950"; new="42×86=3612"},
    @{old="30×19=570"; new="55×46=2530"},
    @{old="50×86=4300"; new="21×14=294"},
    @{old="39×87=3393"; new="23×51=1173"},
    @{old="46×50=2300"; new="95×87=8265"}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $pair.new, 2)
}
